{"js": "// Office.js (Word JavaScript API) script\n// Reproduces: \"adding Viets work back pls continue on seperate sheet\"\n//\n// 1. Left-aligns the existing \"Introduction\" heading.\n// 2. Inserts two body paragraphs (CNN/spectrogram intro + dataset/goal) after it.\n// 3. Promotes the old \"Literature Review\" Heading2 to a left-aligned Heading1,\n//    inserted right before the (now repurposed) heading paragraph.\n// 4. Turns the old \"Literature Review\" paragraph into a new \"Gabor\n//    Transformation\" Heading2, followed by its body paragraphs.\n// 5. Adds a \"Convolutional Neural Network\" Heading2 with its body paragraphs,\n//    a \"....\" placeholder, an \"On going...\" placeholder, and a trailing blank\n//    paragraph, all before the \"(Viet Nguyen)\" signature line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the two anchor paragraphs by their (unique) text.\nlet introPara = null;\nlet literatureReviewPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === \"Introduction\" && introPara === null) {\n    introPara = paragraphs.items[i];\n  } else if (text === \"Literature Review\" && literatureReviewPara === null) {\n    literatureReviewPara = paragraphs.items[i];\n  }\n}\n\nif (!introPara || !literatureReviewPara) {\n  throw new Error(\"Could not locate Introduction / Literature Review headings.\");\n}\n\n// 1. \"Introduction\" heading becomes explicitly left-aligned.\nintroPara.alignment = Word.Alignment.left;\n\n// 2. Two new body paragraphs right after \"Introduction\".\nconst FIRST_LINE_INDENT_PT = 18; // 360 twips\nconst LEFT_INDENT_PT = 14.4;     // 288 twips\n\nconst para1 = introPara.insertParagraph(\n  \"Convolutional neural network (CNN) is a deep learning algorithm used to process the data of image. It is commonly used in computer vision as a classification technique to distinguish different objects. On the other hand, spectrogram is a representation method used to present three-dimension measured signals in two-dimensional diagram. \",\n  Word.InsertLocation.after\n);\npara1.style = \"Normal\";\npara1.alignment = Word.Alignment.left;\npara1.firstLineIndent = FIRST_LINE_INDENT_PT;\n\nconst para2 = para1.insertParagraph(\n  \"Based on the dataset provided by Professor Pech in the module Computational Intelligence at Frankfurt University of Applied Sciences (FRA-UAS), the goal of this project is to classify the reflected signals of different objects using CNN and spectrogram.\",\n  Word.InsertLocation.after\n);\npara2.style = \"Normal\";\npara2.alignment = Word.Alignment.left;\npara2.firstLineIndent = FIRST_LINE_INDENT_PT;\n\n// 3. New \"Literature Review\" Heading1 (left aligned), placed right before the\n//    old \"Literature Review\" paragraph (which becomes \"Gabor Transformation\").\nconst newLitReview = literatureReviewPara.insertParagraph(\"Literature Review\", Word.InsertLocation.before);\nnewLitReview.style = \"Heading 1\";\nnewLitReview.alignment = Word.Alignment.left;\n\n// 4. Repurpose the old \"Literature Review\" Heading2 paragraph's text.\nconst gaborHeading = literatureReviewPara;\ngaborHeading.insertText(\"Gabor Transformation\", Word.InsertLocation.replace);\n// Heading2 style already carries jc=left, so no explicit alignment needed.\n\nconst gaborP1 = gaborHeading.insertParagraph(\n  \"Based on the provided dataset which is the set of analog signals in time domain, Gabor transform is used to convert them to time-frequency representation. Basically, Gabor transform filters the signals with a Gaussian window and Fourier Transform will be then applied to the filtered signals. The following formula is the applied filter as discussed:\",\n  Word.InsertLocation.after\n);\ngaborP1.style = \"Normal\";\ngaborP1.alignment = Word.Alignment.left;\ngaborP1.firstLineIndent = FIRST_LINE_INDENT_PT;\n\nconst gaborFormula = gaborP1.insertParagraph(\n  \"Em chiu thua :'( Giup em cho nay voi, hong hieu gi :v\",\n  Word.InsertLocation.after\n);\ngaborFormula.style = \"Normal\";\ngaborFormula.alignment = Word.Alignment.left;\ngaborFormula.firstLineIndent = FIRST_LINE_INDENT_PT;\n\nconst gaborP2 = gaborFormula.insertParagraph(\n  \"As the time increases, the signal dataset is acquired with the corresponding time from the window length until it reaches the end of the window. The whole process will generate the spectrogram of the signals to be used later as the training set, and also to test the model accuracy.\",\n  Word.InsertLocation.after\n);\ngaborP2.style = \"Normal\";\ngaborP2.alignment = Word.Alignment.left;\ngaborP2.firstLineIndent = FIRST_LINE_INDENT_PT;\n\n// 5. \"Convolutional Neural Network\" Heading2 and its body paragraphs.\nconst cnnHeading = gaborP2.insertParagraph(\"Convolutional Neural Network\", Word.InsertLocation.after);\ncnnHeading.style = \"Heading 2\";\n\nconst cnnP1 = cnnHeading.insertParagraph(\n  \"The concept of neural network or artificial neural network is commonly known as a combination of different layers connected to each other to make decisions based on different types of input. Biologically speaking, the neural network is a technique that mimics approximately how a brain functions. Each layer contains various nodes acts as a system of neurons that can interconnect between layers. Besides, dependent on the importance of each specific neuron, or node, a factor called weight is introduced to bias for the purpose of the system. These layers are commonly known as the hidden layer.\",\n  Word.InsertLocation.after\n);\ncnnP1.style = \"Normal\";\ncnnP1.alignment = Word.Alignment.left;\ncnnP1.leftIndent = LEFT_INDENT_PT;\n\nconst cnnP2 = cnnP1.insertParagraph(\"....\", Word.InsertLocation.after);\ncnnP2.style = \"Normal\";\ncnnP2.alignment = Word.Alignment.left;\ncnnP2.leftIndent = LEFT_INDENT_PT;\n\nconst cnnP3 = cnnP2.insertParagraph(\"On going...\", Word.InsertLocation.after);\ncnnP3.style = \"Normal\";\ncnnP3.alignment = Word.Alignment.left;\ncnnP3.leftIndent = LEFT_INDENT_PT;\n\n// Trailing blank paragraph before \"(Viet Nguyen)\".\nconst blank = cnnP3.insertParagraph(\"\", Word.InsertLocation.after);\nblank.style = \"Normal\";\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Reproduces: \"adding Viets work back pls continue on seperate sheet\"\n#\n# 1. Left-aligns the existing \"Introduction\" heading.\n# 2. Inserts two body paragraphs (CNN/spectrogram intro + dataset/goal) after it.\n# 3. Promotes the old \"Literature Review\" Heading2 to a left-aligned Heading1,\n#    inserted right before the (now repurposed) heading paragraph.\n# 4. Turns the old \"Literature Review\" paragraph into a new \"Gabor\n#    Transformation\" Heading2, followed by its body paragraphs.\n# 5. Adds a \"Convolutional Neural Network\" Heading2 with its body paragraphs,\n#    a \"....\" placeholder, an \"On going...\" placeholder, and a trailing blank\n#    paragraph, all before the \"(Viet Nguyen)\" signature line.\n#\n# NOTE: paragraph COM objects/ranges are NOT held onto across mutations here \u2014\n# every step re-resolves paragraphs via `$d.Paragraphs.Item(<1-based index>)`\n# (or a text search) because inserting a paragraph shifts later indices.\n\n$d = $word.ActiveDocument\n\n$wdAlignParagraphLeft = 0\n$wdCollapseStart = 1\n$wdCollapseEnd = 0\n\nfunction Find-ParaIndexByText($doc, $text) {\n    $idx = 1\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.Trim() -eq $text) {\n            return $idx\n        }\n        $idx++\n    }\n    return -1\n}\n\n# Inserts a new paragraph right AFTER the paragraph currently at 1-based\n# index $idx, fills it in, and returns the new paragraph's 1-based index\n# (always $idx + 1).\nfunction Insert-ParaAfterIdx($doc, $idx, $text, $style, $alignLeft, $firstLineIndent, $leftIndent) {\n    $para = $doc.Paragraphs.Item($idx)\n    $r = $para.Range.Duplicate\n    $r.Collapse($wdCollapseEnd)\n    $r.InsertParagraphAfter()\n\n    $newIdx = $idx + 1\n    $newPara = $doc.Paragraphs.Item($newIdx)\n    $newPara.Range.Text = $text\n    if ($style) { $newPara.Style = $style }\n    if ($alignLeft) { $newPara.Alignment = $wdAlignParagraphLeft }\n    if ($firstLineIndent) { $newPara.Format.FirstLineIndent = $firstLineIndent }\n    if ($leftIndent) { $newPara.Format.LeftIndent = $leftIndent }\n    return $newIdx\n}\n\n# Inserts a new (blank) paragraph right BEFORE the paragraph currently at\n# 1-based index $idx. The new blank paragraph takes over index $idx and the\n# original paragraph's content shifts to $idx + 1. Returns $idx (the new\n# paragraph's index) when $text is supplied (fills it in); otherwise leaves\n# it blank.\nfunction Insert-ParaBeforeIdx($doc, $idx, $text, $style, $alignLeft, $firstLineIndent, $leftIndent) {\n    $para = $doc.Paragraphs.Item($idx)\n    $r = $para.Range.Duplicate\n    $r.Collapse($wdCollapseStart)\n    $r.InsertParagraphBefore()\n\n    $newIdx = $idx\n    if ($text -ne $null) {\n        $newPara = $doc.Paragraphs.Item($newIdx)\n        $newPara.Range.Text = $text\n        if ($style) { $newPara.Style = $style }\n        if ($alignLeft) { $newPara.Alignment = $wdAlignParagraphLeft }\n        if ($firstLineIndent) { $newPara.Format.FirstLineIndent = $firstLineIndent }\n        if ($leftIndent) { $newPara.Format.LeftIndent = $leftIndent }\n    }\n    return $newIdx\n}\n\n$FIRST_LINE_INDENT_PT = 18    # 360 twips\n$LEFT_INDENT_PT = 14.4        # 288 twips\n\n# 1. \"Introduction\" heading becomes explicitly left-aligned.\n$introIdx = Find-ParaIndexByText $d \"Introduction\"\n$d.Paragraphs.Item($introIdx).Alignment = $wdAlignParagraphLeft\n\n# 2. Two new body paragraphs right after \"Introduction\".\n$p1Idx = Insert-ParaAfterIdx $d $introIdx `\n    \"Convolutional neural network (CNN) is a deep learning algorithm used to process the data of image. It is commonly used in computer vision as a classification technique to distinguish different objects. On the other hand, spectrogram is a representation method used to present three-dimension measured signals in two-dimensional diagram. \" `\n    \"Normal\" $true $FIRST_LINE_INDENT_PT $null\n\n$p2Idx = Insert-ParaAfterIdx $d $p1Idx `\n    \"Based on the dataset provided by Professor Pech in the module Computational Intelligence at Frankfurt University of Applied Sciences (FRA-UAS), the goal of this project is to classify the reflected signals of different objects using CNN and spectrogram.\" `\n    \"Normal\" $true $FIRST_LINE_INDENT_PT $null\n\n# 3. New \"Literature Review\" Heading1 (left aligned), inserted right before\n#    the old \"Literature Review\" paragraph (which becomes \"Gabor Transformation\").\n$oldLrIdx = Find-ParaIndexByText $d \"Literature Review\"\n$newLrIdx = Insert-ParaBeforeIdx $d $oldLrIdx \"Literature Review\" \"Heading 1\" $true $null $null\n$gaborHeadingIdx = $newLrIdx + 1   # the old \"Literature Review\" paragraph, shifted down by one\n\n# 4. Repurpose the old \"Literature Review\" Heading2 paragraph's text.\n$d.Paragraphs.Item($gaborHeadingIdx).Range.Text = \"Gabor Transformation\"\n# Heading2 style already carries jc=left, so no explicit alignment needed.\n\n$gaborP1Idx = Insert-ParaAfterIdx $d $gaborHeadingIdx `\n    \"Based on the provided dataset which is the set of analog signals in time domain, Gabor transform is used to convert them to time-frequency representation. Basically, Gabor transform filters the signals with a Gaussian window and Fourier Transform will be then applied to the filtered signals. The following formula is the applied filter as discussed:\" `\n    \"Normal\" $true $FIRST_LINE_INDENT_PT $null\n\n$gaborFormulaIdx = Insert-ParaAfterIdx $d $gaborP1Idx `\n    \"Em chiu thua :'( Giup em cho nay voi, hong hieu gi :v\" `\n    \"Normal\" $true $FIRST_LINE_INDENT_PT $null\n\n$gaborP2Idx = Insert-ParaAfterIdx $d $gaborFormulaIdx `\n    \"As the time increases, the signal dataset is acquired with the corresponding time from the window length until it reaches the end of the window. The whole process will generate the spectrogram of the signals to be used later as the training set, and also to test the model accuracy.\" `\n    \"Normal\" $true $FIRST_LINE_INDENT_PT $null\n\n# 5. \"Convolutional Neural Network\" Heading2 and its body paragraphs.\n$cnnHeadingIdx = Insert-ParaAfterIdx $d $gaborP2Idx \"Convolutional Neural Network\" \"Heading 2\" $false $null $null\n\n$cnnP1Idx = Insert-ParaAfterIdx $d $cnnHeadingIdx `\n    \"The concept of neural network or artificial neural network is commonly known as a combination of different layers connected to each other to make decisions based on different types of input. Biologically speaking, the neural network is a technique that mimics approximately how a brain functions. Each layer contains various nodes acts as a system of neurons that can interconnect between layers. Besides, dependent on the importance of each specific neuron, or node, a factor called weight is introduced to bias for the purpose of the system. These layers are commonly known as the hidden layer.\" `\n    \"Normal\" $true $null $LEFT_INDENT_PT\n\n$cnnP2Idx = Insert-ParaAfterIdx $d $cnnP1Idx \"....\" \"Normal\" $true $null $LEFT_INDENT_PT\n\n$cnnP3Idx = Insert-ParaAfterIdx $d $cnnP2Idx \"On going...\" \"Normal\" $true $null $LEFT_INDENT_PT\n\n# Trailing blank paragraph before \"(Viet Nguyen)\".\n$vietIdx = Find-ParaIndexByText $d \"(Viet Nguyen)\"\nInsert-ParaBeforeIdx $d $vietIdx $null $null $false $null $null | Out-Null\n\n\"ok\"\n"}
